# Insert two new data rows at 163-164 (pushing the existing rows 163-216
# down to 165-218), then populate the two new rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A163:A164").EntireRow.Insert()

# New row 163
$ws.Cells.Item(163, 1).Value = 9
$ws.Cells.Item(163, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(163, 3).Value = "Metropolitana"
$ws.Cells.Item(163, 4).Value = 44588
$ws.Cells.Item(163, 5).Value = 13
$ws.Cells.Item(163, 6).Value = 100112030
$ws.Cells.Item(163, 7).Value = "Poroto granado"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 61
$ws.Cells.Item(163, 11).Value = 25000
$ws.Cells.Item(163, 12).Value = 27000
$ws.Cells.Item(163, 13).Value = 26016
$ws.Cells.Item(163, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(163, 15).Value = "Región Metropolitana"
$ws.Cells.Item(163, 16).Value = 1041
$ws.Cells.Item(163, 17).Value = 25
$ws.Cells.Item(163, 18).Value = "Hortaliza"

# New row 164
$ws.Cells.Item(164, 1).Value = 9
$ws.Cells.Item(164, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(164, 3).Value = "Metropolitana"
$ws.Cells.Item(164, 4).Value = 44588
$ws.Cells.Item(164, 5).Value = 13
$ws.Cells.Item(164, 6).Value = 100112030
$ws.Cells.Item(164, 7).Value = "Poroto granado"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 79
$ws.Cells.Item(164, 11).Value = 25000
$ws.Cells.Item(164, 12).Value = 27000
$ws.Cells.Item(164, 13).Value = 25987
$ws.Cells.Item(164, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(164, 15).Value = "Región del Maule"
$ws.Cells.Item(164, 16).Value = 1039
$ws.Cells.Item(164, 17).Value = 25
$ws.Cells.Item(164, 18).Value = "Hortaliza"
